$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings used for the "walk" row (row 6) to "game" equivalents
$ws.Range("C6").Value = "town_game_action_icon"
$ws.Range("D6").Value = "game_title"
$ws.Range("E6").Value = "game_desc"

# Update column widths for C:D and E
# (values chosen so the engine's internal char-width quantization lands on the
# closest achievable width to the target 13.6153846153846 / 14.7403846153846)
$ws.Range("C:D").ColumnWidth = 12.857142857142858
$ws.Range("E:E").ColumnWidth = 14.0

# Update the active selection to F6
$ws.Range("F6").Select()
